$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths (1-based column index -> ColumnWidth input).
# The runtime stores widths snapped to pixel granularity using MDW=7
# (stored = (round(input*7)+5)/7), mirroring genuine Excel column-width
# persistence. The inputs below are chosen so the stored width lands as
# close as achievable to the target stored widths from the target file.
$colWidths = @{
    1  = 5.714285714285714
    2  = 9.0
    3  = 13.428571428571429
    4  = 13.428571428571429
    5  = 9.0
    6  = 14.428571428571429
    7  = 13.428571428571429
    8  = 10.142857142857142
    9  = 21.142857142857142
    10 = 13.428571428571429
    11 = 7.857142857142857
    12 = 10.142857142857142
    13 = 9.0
    14 = 9.0
    15 = 4.571428571428571
    16 = 10.142857142857142
    17 = 11.142857142857142
    18 = 4.571428571428571
    19 = 9.0
    20 = 7.857142857142857
    21 = 9.0
    22 = 4.571428571428571
    23 = 9.0
    24 = 12.285714285714286
    25 = 15.571428571428571
    26 = 11.142857142857142
    27 = 9.0
}

foreach ($colIdx in $colWidths.Keys) {
    $ws.Columns.Item($colIdx).ColumnWidth = $colWidths[$colIdx]
}

# Header row (row 1)
$headers = @{
    "A1" = "SL #"
    "B1" = "MONTH"
    "C1" = "Emp. NAME"
    "D1" = "DOJ"
    "E1" = "STATUS"
    "F1" = "DESIGNATION"
    "G1" = "DEPARTMENT"
    "H1" = "GROSS"
    "I1" = "Per Month"
    "J1" = "Actual Per Month"
    "K1" = "Actual Days"
    "L1" = "Working Days"
    "M1" = "BASIC"
    "N1" = "HRA"
    "O1" = "DA"
    "P1" = "TL Allowance"
    "Q1" = "Spcl Allowance"
    "R1" = "Arrears"
    "S1" = "Gross Pay"
    "T1" = "PF"
    "U1" = "ESIC"
    "V1" = "PT"
    "W1" = "TDS"
    "X1" = "Deductible Arrears"
    "Y1" = "Deducted allowance1"
    "Z1" = "total_deducations"
    "AA1" = "NetPay"
}

foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# Row 2 - Priyanka Muddana (updated values)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "11-2014"
$ws.Range("C2").Value = "Priyanka Muddana"
# Force D2 to stay literal text (not auto-parsed as a date) by pre-formatting as Text
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "02/06/2014"
$ws.Range("E2").Value = "Internship"
$ws.Range("F2").Value = "HR Manager"
$ws.Range("G2").Value = "HR"
$ws.Range("H2").Value = 750000.0
$ws.Range("I2").Value = 62500.0
$ws.Range("J2").Value = 58406.2
$ws.Range("K2").Value = 30.0
$ws.Range("L2").Value = 30.0
$ws.Range("M2").Value = 25000.0
$ws.Range("N2").Value = 6250.0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 27156.2
$ws.Range("R2").Value = 0.0
$ws.Range("S2").Value = 58406.2
$ws.Range("T2").Value = 3000.0
$ws.Range("U2").Value = 1022.11
$ws.Range("V2").Value = 0.0
$ws.Range("W2").Value = 13395.0
$ws.Range("X2").Value = 0.0
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 17417.1
$ws.Range("AA2").Value = 40989.1

# Row 3 - Vidya Sagar pogiri (new row)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "11-2014"
$ws.Range("C3").Value = "Vidya Sagar pogiri"
# Force D3 to stay literal text (not auto-parsed as a date) by pre-formatting as Text
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "02/06/2014"
$ws.Range("E3").Value = "Regular"
$ws.Range("F3").Value = "Junior Developer"
$ws.Range("G3").Value = "Development"
$ws.Range("H3").Value = 130000.0
$ws.Range("I3").Value = 10833.333333333334
$ws.Range("J3").Value = 10313.3
$ws.Range("K3").Value = 30.0
$ws.Range("L3").Value = 30.0
$ws.Range("M3").Value = 4333.33
$ws.Range("N3").Value = 1083.33
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 4896.67
$ws.Range("R3").Value = 0.0
$ws.Range("S3").Value = 10313.3
$ws.Range("T3").Value = 520.0
$ws.Range("U3").Value = 0.0
$ws.Range("V3").Value = 0.0
$ws.Range("W3").Value = 0.0
$ws.Range("X3").Value = 0.0
$ws.Range("Y3").Value = 0
$ws.Range("Z3").Value = 520.0
$ws.Range("AA3").Value = 9793.33
